# Update the gyroscope data block (rows 2:31, columns A:C) with the new
# window of readings. The window slid forward by 11 rows (old rows 13-21
# are now rows 2-10) and grew to include 21 additional new rows, extending
# the sheet from A1:C21 to A1:C31.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
        @(0.1481355747003131, -0.5275314320756641, -1.595921289858394),
        @(0.8970806586181024, -2.253284974872135, -2.632813992924737),
        @(3.44800701316113, -1.11263772954487, -4.982394183493419),
        @(-0.06218727720976736, 0.8519693432053494, -3.113853278584423),
        @(-5.216508371043581, -1.343549091778454, 3.495817017180781),
        @(-5.247863911833437, 3.879994811812087, 4.198587553663369),
        @(-3.767358408548382, 6.815083000828341, 0.4525070901939874),
        @(1.195290473119145, -1.12244814540705, -3.771863303259405),
        @(4.371577180492927, -3.279605823973685, -4.797551259944596),
        @(0.8983272033211789, 1.015922405332796, -5.503626271692188),
        @(-3.768744276456209, 1.106249593939454, 3.618037942816487),
        @(-5.589900377533165, 3.503953431913461, 6.05607688553742),
        @(-2.11142427696607, 9.510962799581424, -0.04490738020751496),
        @(0.3853269297416066, -2.830002902065919, -0.4395293002041412),
        @(7.810094511321666, -5.127505102706861, -1.701906724125914),
        @(4.436884702500196, -6.79071483312473, -0.1620084203350239),
        @(-0.8537680622172118, -3.271142317050352, 3.793036780432373),
        @(-7.225753661849682, 1.326318728362095, 8.623776784117876),
        @(3.725271386937927, 3.283059017820457, -3.099974366226723),
        @(2.260886770268454, -3.733322253401921, 0.4708228362480379),
        @(10.02548202544597, 5.10485092133122, 0.5539246020516877),
        @(-2.205630293691351, -1.242692031161299, 1.271664639886122),
        @(-2.934371616203759, 0.8867856582421496, -0.5782981265278011),
        @(-3.666905725004196, -0.9611553521678435, -2.932546598124879),
        @(0.07597643337871512, 15.43764277153591, -3.051599523783981),
        @(0.9658543077439407, -2.044685548512497, -3.529194889268377),
        @(4.105489831320192, 13.38752136726634, -6.598481021626331),
        @(-1.887021554078156, 0.322081098381812, 1.926029883130062),
        @(-4.947177827046182, 3.337713633454262, 1.887203427509456),
        @(-5.072363108240451, -0.134192296503722, -3.107403220930567)
)

$startRow = 2

for ($i = 0; $i -lt $data.Length; $i++) {
    $rowValues = $data[$i]
    $targetRow = $startRow + $i
    for ($j = 0; $j -lt $rowValues.Length; $j++) {
        $ws.Cells.Item($targetRow, $j + 1).Value = $rowValues[$j]
    }
}
